$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 93.62780766666667
$ws.Range("H2").Value = 280.883423
$ws.Range("I2").Value = 0.3228593149748609
$ws.Range("J2").Value = 0.3228593149748609
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 968.824382191638
$ws.Range("R2").Value = 8719.419439724743
$ws.Range("S2").Value = 0.07425614042590413
$ws.Range("T2").Value = 0.07425614042590414

# Row 3
$ws.Range("G3").Value = 93.62780766666667
$ws.Range("H3").Value = 280.883423
$ws.Range("I3").Value = 0.3228593149748609
$ws.Range("J3").Value = 0.3228593149748609
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("Q3").Value = 2862.189190941811
$ws.Range("R3").Value = 25759.7027184763
$ws.Range("S3").Value = 0.2193742502715416
$ws.Range("T3").Value = 0.2193742502715416

# Row 4
$ws.Range("G4").Value = 93.62780766666667
$ws.Range("H4").Value = 280.883423
$ws.Range("I4").Value = 0.3228593149748609
$ws.Range("J4").Value = 0.3228593149748609
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 381.3515534577164
$ws.Range("R4").Value = 3432.163981119448
$ws.Range("S4").Value = 0.02922892427741512
$ws.Range("T4").Value = 0.02922892427741513

# Row 5
$ws.Range("G5").Value = 66.39541
$ws.Range("I5").Value = 0.228953097635189
$ws.Range("J5").Value = 0.228953097635189
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 687.0340519199366
$ws.Range("R5").Value = 6183.30646727943
$ws.Range("S5").Value = 0.05265814731183491
$ws.Range("T5").Value = 0.05265814731183492

# Row 6
$ws.Range("G6").Value = 66.39541
$ws.Range("I6").Value = 0.228953097635189
$ws.Range("J6").Value = 0.228953097635189
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("Q6").Value = 2029.698543265223
$ws.Range("S6").Value = 0.155567492748281
$ws.Range("T6").Value = 0.155567492748281

# Row 7
$ws.Range("G7").Value = 66.39541
$ws.Range("I7").Value = 0.228953097635189
$ws.Range("J7").Value = 0.228953097635189
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 270.4324001273866
$ws.Range("R7").Value = 2433.89160114648
$ws.Range("S7").Value = 0.02072745757507303
$ws.Range("T7").Value = 0.02072745757507303

# Row 8
$ws.Range("G8").Value = 129.9724656666667
$ws.Range("H8").Value = 389.917397
$ws.Range("I8").Value = 0.4481875873899502
$ws.Range("J8").Value = 0.4481875873899502
$ws.Range("M8").Value = 10.34761366666667
$ws.Range("N8").Value = 31.042841
$ws.Range("O8").Value = 0.2299953477621856
$ws.Range("P8").Value = 0.2299953477621856
$ws.Range("Q8").Value = 1344.904862022764
$ws.Range("R8").Value = 12104.14375820488
$ws.Range("S8").Value = 0.1030810600244465
$ws.Range("T8").Value = 0.1030810600244466

# Row 9
$ws.Range("G9").Value = 129.9724656666667
$ws.Range("H9").Value = 389.917397
$ws.Range("I9").Value = 0.4481875873899502
$ws.Range("J9").Value = 0.4481875873899502
$ws.Range("O9").Value = 0.6794731949692173
$ws.Range("P9").Value = 0.6794731949692174
$ws.Range("Q9").Value = 3973.240382553893
$ws.Range("R9").Value = 35759.16344298504
$ws.Range("S9").Value = 0.3045314519493948
$ws.Range("T9").Value = 0.3045314519493948

# Row 10
$ws.Range("G10").Value = 129.9724656666667
$ws.Range("H10").Value = 389.917397
$ws.Range("I10").Value = 0.4481875873899502
$ws.Range("J10").Value = 0.4481875873899502
$ws.Range("M10").Value = 4.073058666666666
$ws.Range("N10").Value = 12.219176
$ws.Range("O10").Value = 0.09053145726859702
$ws.Range("P10").Value = 0.09053145726859703
$ws.Range("Q10").Value = 529.3854777116524
$ws.Range("R10").Value = 4764.469299404872
$ws.Range("S10").Value = 0.04057507541610887
$ws.Range("T10").Value = 0.04057507541610888
